$d = $word.ActiveDocument

# --- 1. Remove the three obsolete to-do items -----------------------------
# "Create a quick title screen.", "Create a Game Scene" and
# "Add player collision so that player can die." are dropped; the list now
# starts with "Add asteroid collision ..." (formerly the 4th bullet).
$pStart = $d.Paragraphs(3)
$pEnd   = $d.Paragraphs(5)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# --- 2. Move the _GoBack bookmark from the last bullet to the new first --
# bullet ("Add asteroid collision so that player can collide with
# player."), placing it right after that run's text.
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$targetPara = $d.Paragraphs(3)
$insertPos = $targetPara.Range.End - 1

# Word's bookmark engine mishandles a zero-length range that sits exactly
# at "end of paragraph text" (immediately before the paragraph mark), so a
# harmless placeholder character is inserted first to give the bookmark a
# safe, unambiguous anchor; it is removed again right after.
$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanup = $d.Range($insertPos, $insertPos + 1)
$cleanup.Delete()
